# Se agrega validación de fecha para Ecuaquimica
#
# Ecuaquimica ("ECUATORIANA DE PROD. QUIM. S.A") is the client listed in
# row 2 of the "Base Clientes carga manual" sheet. Its distributor number
# (Num_Distri) is corrected from 61610097 to 61610107.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base Clientes carga manual")
$ws.Activate()

# Correct the distributor number for Ecuaquimica (A2, table column Num_Distri).
$ws.Range("A2").Value = 61610107

# Leave the selection where the user ended up after making the edit.
$ws.Range("A3").Select()

# Reposition/resize the Excel window to match the state captured on save.
$win = $excel.ActiveWindow
$win.Left = 26280
$win.Top = 1170
$win.Width = 21840
$win.Height = 13020
